$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 176
$ws.Range("J2").Value = 96.86667
$ws.Range("L2").Value = 96.86667
$ws.Range("N2").Value = -322.86667
$ws.Range("H9").Value = 92
$ws.Range("I9").Value = 86.9375
$ws.Range("J9").Value = 102.125
$ws.Range("K9").Value = 86.9375
$ws.Range("L9").Value = 102.125
$ws.Range("M9").Value = 82.0625
$ws.Range("N9").Value = -440.125
$ws.Range("H19").Value = 1035.5333
$ws.Range("I19").Value = 407.25
$ws.Range("K19").Value = 407.25
$ws.Range("M19").Value = -232.25
$ws.Range("H28").Value = 1795.9615
$ws.Range("I28").Value = 1668.2273
$ws.Range("K28").Value = 1668.2273
$ws.Range("M28").Value = -1183.2273
$ws.Range("H33").Value = 72143220
$ws.Range("I33").Value = 833584.25
$ws.Range("J33").Value = 500001000
$ws.Range("K33").Value = 833584.25
$ws.Range("L33").Value = 500001000
$ws.Range("M33").Value = -833355.25
$ws.Range("N33").Value = -500001458
$ws.Range("H58").Value = 4256
$ws.Range("I58").Value = 413
$ws.Range("J58").Value = 35000
$ws.Range("K58").Value = 1239
$ws.Range("L58").Value = 105000
$ws.Range("M58").Value = -1089
$ws.Range("N58").Value = -105300
$ws.Range("H62").Value = 9442.5
$ws.Range("I62").Value = 9442.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 9442.5
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -8818.5
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 9442.5
$ws.Range("I65").Value = 9442.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 47212.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -44092.5
$ws.Range("N65").ClearContents()
$ws.Range("H94").Value = 2679.9
$ws.Range("I94").Value = 2679.9
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2679.9
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -2228.9
$ws.Range("N94").ClearContents()
$ws.Range("H106").Value = 55821.09
$ws.Range("I106").Value = 71628.875
$ws.Range("K106").Value = 71628.875
$ws.Range("M106").Value = -70997.875
$ws.Range("H112").Value = 1256.629
$ws.Range("J112").Value = 1253.5737
$ws.Range("L112").Value = 3760.7211
$ws.Range("N112").Value = -5976.7211
$ws.Range("H113").Value = 5866.7144
$ws.Range("I113").Value = 7387.5386
$ws.Range("J113").Value = 5519.86
$ws.Range("K113").Value = 7387.5386
$ws.Range("L113").Value = 5519.86
$ws.Range("M113").Value = -4133.5386
$ws.Range("N113").Value = -12027.86
$ws.Range("H116").Value = 2375.8462
$ws.Range("I116").Value = 2485
$ws.Range("J116").Value = 2201.2
$ws.Range("K116").Value = 2485
$ws.Range("L116").Value = 2201.2
$ws.Range("M116").Value = 957
$ws.Range("N116").Value = -9085.200000000001
$ws.Range("H118").Value = 646.75
$ws.Range("J118").Value = 499.5
$ws.Range("L118").Value = 1498.5
$ws.Range("N118").Value = -4812.5
$ws.Range("H127").Value = 8499
$ws.Range("I127").Value = 8499
$ws.Range("K127").Value = 25497
$ws.Range("M127").Value = -20537
$ws.Range("H132").Value = 2016.186
$ws.Range("I132").Value = 1937.9459
$ws.Range("J132").Value = 2498.6667
$ws.Range("K132").Value = 5813.8377
$ws.Range("L132").Value = 7496.000100000001
$ws.Range("M132").Value = -3283.8377
$ws.Range("N132").Value = -12556.0001
$ws.Range("H135").Value = 2219.4614
$ws.Range("I135").Value = 2205.7083
$ws.Range("K135").Value = 19851.3747
$ws.Range("M135").Value = -17316.3747
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 2231.2693
$ws.Range("I137").Value = 2586.5
$ws.Range("J137").Value = 1662.9
$ws.Range("K137").Value = 7759.5
$ws.Range("L137").Value = 4988.700000000001
$ws.Range("M137").Value = -5209.5
$ws.Range("N137").Value = -10088.7
$ws.Range("H138").Value = 15875385
$ws.Range("I138").Value = 83334900
$ws.Range("J138").Value = 2559.843
$ws.Range("K138").Value = 250004700
$ws.Range("L138").Value = 7679.529
$ws.Range("M138").Value = -249999560
$ws.Range("N138").Value = -17959.529

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9544.846
$ws.Range("I2").Value = 1158.7
$ws.Range("J2").Value = 37498.668
$ws.Range("K2").Value = 1158.7
$ws.Range("L2").Value = 37498.668
$ws.Range("M2").Value = -1045.7
$ws.Range("N2").Value = -37724.668
$ws.Range("H32").Value = 3111.25
$ws.Range("I32").Value = 1919.7593
$ws.Range("K32").Value = 1919.7593
$ws.Range("M32").Value = -1632.7593
$ws.Range("H45").Value = 5332
$ws.Range("I45").Value = 5332
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 5332
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -4955
$ws.Range("N45").ClearContents()
$ws.Range("H48").Value = 118992.5
$ws.Range("J48").Value = 118992.5
$ws.Range("L48").Value = 118992.5
$ws.Range("N48").Value = -119760.5
$ws.Range("H61").Value = 2323.7693
$ws.Range("I61").Value = 1871
$ws.Range("J61").Value = 3833
$ws.Range("K61").Value = 1871
$ws.Range("L61").Value = 3833
$ws.Range("M61").Value = -1659
$ws.Range("N61").Value = -4257
$ws.Range("H63").Value = 3018.182
$ws.Range("I63").Value = 3133.4443
$ws.Range("J63").Value = 2499.5
$ws.Range("K63").Value = 3133.4443
$ws.Range("L63").Value = 2499.5
$ws.Range("M63").Value = -2447.4443
$ws.Range("N63").Value = -3871.5
$ws.Range("H66").Value = 3018.182
$ws.Range("I66").Value = 3133.4443
$ws.Range("J66").Value = 2499.5
$ws.Range("K66").Value = 15667.2215
$ws.Range("L66").Value = 12497.5
$ws.Range("M66").Value = -12235.2215
$ws.Range("N66").Value = -19361.5
$ws.Range("H74").Value = 3014.4773
$ws.Range("I74").Value = 2572.6897
$ws.Range("J74").Value = 3868.6
$ws.Range("K74").Value = 2572.6897
$ws.Range("L74").Value = 3868.6
$ws.Range("M74").Value = -1698.6897
$ws.Range("N74").Value = -5616.6
$ws.Range("H77").Value = 3014.4773
$ws.Range("I77").Value = 2572.6897
$ws.Range("J77").Value = 3868.6
$ws.Range("K77").Value = 12863.4485
$ws.Range("L77").Value = 19343
$ws.Range("M77").Value = -8495.448499999999
$ws.Range("N77").Value = -28079
$ws.Range("H88").Value = 114443.664
$ws.Range("I88").Value = 333999.66
$ws.Range("K88").Value = 333999.66
$ws.Range("M88").Value = -333593.66
$ws.Range("H91").Value = 114443.664
$ws.Range("I91").Value = 333999.66
$ws.Range("K91").Value = 333999.66
$ws.Range("M91").Value = -332595.66
$ws.Range("H102").Value = 37895.715
$ws.Range("I102").Value = 1712.1666
$ws.Range("K102").Value = 1712.1666
$ws.Range("M102").Value = -90.16660000000002
$ws.Range("H110").Value = 1235.1818
$ws.Range("I110").Value = 1235.1818
$ws.Range("K110").Value = 1235.1818
$ws.Range("M110").Value = 809.8181999999999
$ws.Range("H116").Value = 9544.846
$ws.Range("I116").Value = 1158.7
$ws.Range("J116").Value = 37498.668
$ws.Range("K116").Value = 1158.7
$ws.Range("L116").Value = 37498.668
$ws.Range("M116").Value = 1135.3
$ws.Range("N116").Value = -42086.668
$ws.Range("H122").Value = 5656.815
$ws.Range("I122").Value = 5641.8
$ws.Range("J122").Value = 5699.7144
$ws.Range("K122").Value = 16925.4
$ws.Range("L122").Value = 17099.1432
$ws.Range("M122").Value = -14475.4
$ws.Range("N122").Value = -21999.1432
$ws.Range("H132").Value = 29250.406
$ws.Range("I132").Value = 30105.188
$ws.Range("J132").Value = 23033.818
$ws.Range("K132").Value = 90315.564
$ws.Range("L132").Value = 69101.454
$ws.Range("M132").Value = -87785.564
$ws.Range("N132").Value = -74161.454
$ws.Range("H136").Value = 2323.7693
$ws.Range("I136").Value = 1871
$ws.Range("J136").Value = 3833
$ws.Range("K136").Value = 5613
$ws.Range("L136").Value = 11499
$ws.Range("M136").Value = -3063
$ws.Range("N136").Value = -16599

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9544.846
$ws.Range("I3").Value = 1158.7
$ws.Range("J3").Value = 37498.668
$ws.Range("K3").Value = 1158.7
$ws.Range("L3").Value = 37498.668
$ws.Range("M3").Value = -1044.7
$ws.Range("N3").Value = -37726.668
$ws.Range("H12").Value = 2073.5715
$ws.Range("I12").Value = 834.6667
$ws.Range("J12").Value = 3002.75
$ws.Range("K12").Value = 834.6667
$ws.Range("L12").Value = 3002.75
$ws.Range("M12").Value = -666.6667
$ws.Range("N12").Value = -3338.75
$ws.Range("H20").Value = 2946
$ws.Range("I20").Value = 3187.6
$ws.Range("J20").Value = 2644
$ws.Range("K20").Value = 3187.6
$ws.Range("L20").Value = 2644
$ws.Range("M20").Value = -2940.6
$ws.Range("N20").Value = -3138
$ws.Range("H42").Value = 119947.5
$ws.Range("J42").Value = 119947.5
$ws.Range("L42").Value = 119947.5
$ws.Range("N42").Value = -120603.5
$ws.Range("H43").Value = 171638.33
$ws.Range("J43").Value = 171638.33
$ws.Range("L43").Value = 171638.33
$ws.Range("N43").Value = -172000.33
$ws.Range("H44").Value = 26924
$ws.Range("J44").Value = 26924
$ws.Range("L44").Value = 26924
$ws.Range("N44").Value = -27918
$ws.Range("H76").Value = 52650.6
$ws.Range("J76").Value = 52650.6
$ws.Range("L76").Value = 52650.6
$ws.Range("N76").Value = -53280.6
$ws.Range("H79").Value = 52650.6
$ws.Range("J79").Value = 52650.6
$ws.Range("L79").Value = 52650.6
$ws.Range("N79").Value = -54834.6
$ws.Range("H86").Value = 3708.182
$ws.Range("I86").Value = 3250
$ws.Range("J86").Value = 4090
$ws.Range("K86").Value = 3250
$ws.Range("L86").Value = 4090
$ws.Range("M86").Value = -2127
$ws.Range("N86").Value = -6336
$ws.Range("H89").Value = 3708.182
$ws.Range("I89").Value = 3250
$ws.Range("J89").Value = 4090
$ws.Range("K89").Value = 16250
$ws.Range("L89").Value = 20450
$ws.Range("M89").Value = -10634
$ws.Range("N89").Value = -31682
$ws.Range("H99").Value = 142858800
$ws.Range("I99").Value = 200001150
$ws.Range("J99").Value = 2909.5
$ws.Range("K99").Value = 200001150
$ws.Range("L99").Value = 2909.5
$ws.Range("M99").Value = -199999652
$ws.Range("N99").Value = -5905.5
$ws.Range("H103").Value = 11250
$ws.Range("J103").Value = 11250
$ws.Range("L103").Value = 11250
$ws.Range("N103").Value = -13594
$ws.Range("H105").Value = 4710.222
$ws.Range("J105").Value = 4712.5
$ws.Range("L105").Value = 4712.5
$ws.Range("N105").Value = -8206.5
$ws.Range("H107").Value = 1598.75
$ws.Range("I107").Value = 1598.75
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1598.75
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 321.25
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 2751.9583
$ws.Range("I134").Value = 2807.6191
$ws.Range("J134").Value = 2362.3333
$ws.Range("K134").Value = 8422.8573
$ws.Range("L134").Value = 7086.999899999999
$ws.Range("M134").Value = -5887.8573
$ws.Range("N134").Value = -12156.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 1599.8
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 1599.8
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 1599.8
$ws.Range("N25").Value = -1947.8
$ws.Range("M25").ClearContents()
$ws.Range("H31").Value = 1756.2727
$ws.Range("I31").Value = 1516.1666
$ws.Range("J31").Value = 2044.4
$ws.Range("K31").Value = 1516.1666
$ws.Range("L31").Value = 2044.4
$ws.Range("M31").Value = -1221.1666
$ws.Range("N31").Value = -2634.4
$ws.Range("H34").Value = 1756.2727
$ws.Range("I34").Value = 1516.1666
$ws.Range("J34").Value = 2044.4
$ws.Range("K34").Value = 1516.1666
$ws.Range("L34").Value = 2044.4
$ws.Range("M34").Value = -1314.1666
$ws.Range("N34").Value = -2448.4
$ws.Range("H58").Value = 3595.0833
$ws.Range("I58").Value = 3737.889
$ws.Range("K58").Value = 3737.889
$ws.Range("M58").Value = -3534.889
$ws.Range("H74").Value = 28404
$ws.Range("J74").Value = 28404
$ws.Range("L74").Value = 28404
$ws.Range("N74").Value = -30152
$ws.Range("H77").Value = 28404
$ws.Range("J77").Value = 28404
$ws.Range("L77").Value = 85212
$ws.Range("N77").Value = -93948
$ws.Range("H99").Value = 3835.7222
$ws.Range("I99").Value = 1895.3636
$ws.Range("K99").Value = 1895.3636
$ws.Range("M99").Value = -397.3635999999999
$ws.Range("H107").Value = 1088.7435
$ws.Range("I107").Value = 1111.6666
$ws.Range("K107").Value = 1111.6666
$ws.Range("M107").Value = 808.3334
$ws.Range("H122").Value = 2378.5
$ws.Range("I122").Value = 2307.4443
$ws.Range("J122").Value = 2449.5557
$ws.Range("K122").Value = 6922.3329
$ws.Range("L122").Value = 7348.6671
$ws.Range("M122").Value = -4472.3329
$ws.Range("N122").Value = -12248.6671
$ws.Range("H126").Value = 3835.7222
$ws.Range("I126").Value = 1895.3636
$ws.Range("K126").Value = 5686.0908
$ws.Range("M126").Value = -3216.0908
$ws.Range("H132").Value = 1311
$ws.Range("I132").Value = 1311
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3933
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1403
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 1530.5416
$ws.Range("I134").Value = 1607.0952
$ws.Range("J134").Value = 994.6667
$ws.Range("K134").Value = 4821.2856
$ws.Range("L134").Value = 2984.0001
$ws.Range("M134").Value = -2286.2856
$ws.Range("N134").Value = -8054.0001
$ws.Range("H136").Value = 3595.0833
$ws.Range("I136").Value = 3737.889
$ws.Range("K136").Value = 11213.667
$ws.Range("M136").Value = -8663.667000000001
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 14492969
$ws.Range("J2").Value = 23809756
$ws.Range("L2").Value = 142858536
$ws.Range("N2").Value = -142858762
$ws.Range("H11").Value = 416
$ws.Range("I11").Value = 349.5
$ws.Range("K11").Value = 1048.5
$ws.Range("M11").Value = -908.5
$ws.Range("H17").Value = 3799.3333
$ws.Range("J17").Value = 5500
$ws.Range("L17").Value = 16500
$ws.Range("N17").Value = -16838
$ws.Range("H60").Value = 249
$ws.Range("I60").Value = 249
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 747
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("M60").Value = -496
$ws.Range("H92").Value = 368.4
$ws.Range("J92").Value = 349.5
$ws.Range("L92").Value = 1048.5
$ws.Range("N92").Value = -3544.5
$ws.Range("H94").Value = 10921.25
$ws.Range("J94").Value = 17997.5
$ws.Range("L94").Value = 53992.5
$ws.Range("N94").Value = -55344.5
$ws.Range("H131").Value = 955.5294
$ws.Range("J131").Value = 1895.25
$ws.Range("L131").Value = 5685.75
$ws.Range("N131").Value = -15765.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 29999
$ws.Range("J47").Value = 29999
$ws.Range("L47").Value = 29999
$ws.Range("N47").Value = -31135
$ws.Range("H70").Value = 75199.22
$ws.Range("I70").Value = 99467.164
$ws.Range("J70").Value = 26663.334
$ws.Range("K70").Value = 99467.164
$ws.Range("L70").Value = 26663.334
$ws.Range("M70").Value = -99197.164
$ws.Range("N70").Value = -27203.334
$ws.Range("H73").Value = 75199.22
$ws.Range("I73").Value = 99467.164
$ws.Range("J73").Value = 26663.334
$ws.Range("K73").Value = 99467.164
$ws.Range("L73").Value = 26663.334
$ws.Range("M73").Value = -98531.164
$ws.Range("N73").Value = -28535.334
$ws.Range("H80").Value = 3366.9333
$ws.Range("J80").Value = 3101.9
$ws.Range("L80").Value = 3101.9
$ws.Range("N80").Value = -5097.9
$ws.Range("H83").Value = 3366.9333
$ws.Range("J83").Value = 3101.9
$ws.Range("L83").Value = 15509.5
$ws.Range("N83").Value = -25493.5
$ws.Range("H102").Value = 3579
$ws.Range("I102").Value = 3579
$ws.Range("K102").Value = 3579
$ws.Range("M102").Value = -1957
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1000
$ws.Range("N113").Value = -5340
$ws.Range("M113").ClearContents()
$ws.Range("H126").Value = 4857.25
$ws.Range("I126").Value = 4122.5713
$ws.Range("K126").Value = 12367.7139
$ws.Range("M126").Value = -9897.713899999999
$ws.Range("H132").Value = 1813.6207
$ws.Range("I132").Value = 1813.6207
$ws.Range("K132").Value = 5440.8621
$ws.Range("M132").Value = -2910.8621

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5706.8184
$ws.Range("I7").Value = 3254
$ws.Range("J7").Value = 9999.25
$ws.Range("K7").Value = 3254
$ws.Range("L7").Value = 9999.25
$ws.Range("M7").Value = -3142
$ws.Range("N7").Value = -10223.25
$ws.Range("H22").Value = 114675.445
$ws.Range("J22").Value = 7274.75
$ws.Range("L22").Value = 7274.75
$ws.Range("N22").Value = -7864.75
$ws.Range("H27").Value = 114675.445
$ws.Range("J27").Value = 7274.75
$ws.Range("L27").Value = 7274.75
$ws.Range("N27").Value = -7488.75
$ws.Range("H40").Value = 5048
$ws.Range("I40").Value = 4916
$ws.Range("K40").Value = 4916
$ws.Range("M40").Value = -4780
$ws.Range("H46").Value = 882
$ws.Range("I46").Value = 909.6667
$ws.Range("J46").Value = 799
$ws.Range("K46").Value = 909.6667
$ws.Range("L46").Value = 799
$ws.Range("M46").Value = -721.6667
$ws.Range("N46").Value = -1175
$ws.Range("H61").Value = 1519.5
$ws.Range("I61").Value = 1330.4
$ws.Range("J61").Value = 1834.6666
$ws.Range("K61").Value = 1330.4
$ws.Range("L61").Value = 1834.6666
$ws.Range("M61").Value = -1128.4
$ws.Range("N61").Value = -2238.6666
$ws.Range("H82").Value = 1734.5758
$ws.Range("I82").Value = 986.55554
$ws.Range("K82").Value = 986.55554
$ws.Range("M82").Value = -625.55554
$ws.Range("H85").Value = 1734.5758
$ws.Range("I85").Value = 986.55554
$ws.Range("K85").Value = 986.55554
$ws.Range("M85").Value = 261.44446
$ws.Range("H113").Value = 1519.5
$ws.Range("I113").Value = 1330.4
$ws.Range("J113").Value = 1834.6666
$ws.Range("K113").Value = 1330.4
$ws.Range("L113").Value = 1834.6666
$ws.Range("M113").Value = 839.5999999999999
$ws.Range("N113").Value = -6174.6666
$ws.Range("H122").Value = 5791.95
$ws.Range("I122").Value = 4511.364
$ws.Range("K122").Value = 13534.092
$ws.Range("M122").Value = -11084.092
$ws.Range("H126").Value = 5706.8184
$ws.Range("I126").Value = 3254
$ws.Range("J126").Value = 9999.25
$ws.Range("K126").Value = 9762
$ws.Range("L126").Value = 29997.75
$ws.Range("M126").Value = -7292
$ws.Range("N126").Value = -34937.75
$ws.Range("H132").Value = 3558.0852
$ws.Range("I132").Value = 3070.7097
$ws.Range("K132").Value = 9212.1291
$ws.Range("M132").Value = -6682.1291
$ws.Range("H136").Value = 5440.2666
$ws.Range("I136").Value = 5114.643
$ws.Range("K136").Value = 15343.929
$ws.Range("M136").Value = -12793.929

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 18056.334
$ws.Range("J41").Value = 18073.572
$ws.Range("L41").Value = 18073.572
$ws.Range("N41").Value = -18853.572
$ws.Range("H81").Value = 7241.3335
$ws.Range("I81").Value = 3316.6667
$ws.Range("K81").Value = 6633.3334
$ws.Range("M81").Value = -5572.3334
$ws.Range("H84").Value = 7241.3335
$ws.Range("I84").Value = 3316.6667
$ws.Range("K84").Value = 33166.667
$ws.Range("M84").Value = -27862.667
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H122").Value = 6581.8125
$ws.Range("I122").Value = 6765.952
$ws.Range("J122").Value = 6230.273
$ws.Range("K122").Value = 20297.856
$ws.Range("L122").Value = 18690.819
$ws.Range("M122").Value = -17847.856
$ws.Range("N122").Value = -23590.819
$ws.Range("H123").Value = 64249.5
$ws.Range("J123").Value = 64249.5
$ws.Range("L123").Value = 64249.5
$ws.Range("N123").Value = -74049.5
$ws.Range("H132").Value = 857.44446
$ws.Range("I132").Value = 889.75
$ws.Range("K132").Value = 2669.25
$ws.Range("M132").Value = -139.25
$ws.Range("H136").Value = 6378.161
$ws.Range("I136").Value = 5958.846
$ws.Range("K136").Value = 17876.538
$ws.Range("M136").Value = -15326.538
